# updated samples according to rggen/rggen#8
#
# register_5 gains two additional bit fields (bit_field_2 and bit_field_3),
# splitting what used to be two 4-bit wide fields into four 2-bit wide
# fields. Everything that used to follow register_5 (register_6, register_7,
# register_8, register_9) shifts down by two rows as a result.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows right after the current bit_field_1 row of
# register_5 (row 22), pushing bit_field_2.. (and every register below)
# down by two rows.
$ws.Rows("22:23").Insert()

# The inserted rows don't inherit the bordered "middle of block" styling
# used by columns B:J, so copy it over from the row directly above
# (which already uses that same style) - restricted to B:J only so we
# don't clobber the rest of the (unbounded) row.
$ws.Range("B21:J21").Copy()
$ws.Range("B22:J23").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# register_5.bit_field_0 / bit_field_1 shrink from 4 bits to 2 bits each
$ws.Range("G20").Value = "0:2"
$ws.Range("G21").Value = "2:2"

# New register_5.bit_field_2
$ws.Range("F22").Value = "bit_field_2"
$ws.Range("G22").Value = "4:2"
$ws.Range("H22").Value = "rws"
$ws.Range("I22").Value = 0
$ws.Range("J22").Value = ""

# New register_5.bit_field_3
$ws.Range("F23").Value = "bit_field_3"
$ws.Range("G23").Value = "6:2"
$ws.Range("H23").Value = "rws"
$ws.Range("I23").Value = 0
$ws.Range("J23").Value = "register_3.bit_field_2"

# The remaining register_5 bit fields (previously bit_field_2..bit_field_7,
# now shifted down to rows 24-29) are renumbered to bit_field_4..bit_field_9.
$ws.Range("F24").Value = "bit_field_4"
$ws.Range("F25").Value = "bit_field_5"
$ws.Range("F26").Value = "bit_field_6"
$ws.Range("F27").Value = "bit_field_7"
$ws.Range("F28").Value = "bit_field_8"
$ws.Range("F29").Value = "bit_field_9"
